$wb = $excel.ActiveWorkbook

$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")
$assets    = $wb.Worksheets.Item("Assets")

# --- Settings sheet: update / add the AdminInvoiceHandler related rows ---
# Row 2: Orchestrator queue name asset
$settings.Range("B2").Value = "AdminInvoicesS2E"
# Row 4: business process name value
$settings.Range("B4").Value = "AdminInvoiceHandler"
# Row 6 (new): Invoice path setting
$settings.Range("A6").Value = "InvoicePath"
$settings.Range("B6").Value = "Invoices"
# Row 8 (new): Business units setting
$settings.Range("A8").Value = "BusinessUnits"
$settings.Range("B8").Value = "AFCL, MKTG, DDM, IRE, SEC, SALES, DIR, PPL, PMI, ISI, ISI-ICT, terze parti-partite iva"
# Re-assert the Name / Description columns for rows 2 and 4 (unchanged text, keeps things tidy)
$settings.Range("A2").Value = "OrchestratorQueueName"
$settings.Range("C2").Value = "Orchestrator queue Name. The value must match with the queue name defined on Orchestrator."
$settings.Range("C4").Value = "Logging field which allows grouping of log data of two or more subprocesses under the same business process name"

# Row heights: row 4 shrinks slightly on both Settings and Constants
$settings.Rows.Item(4).RowHeight = 29
$constants.Rows.Item(2).RowHeight = 29

# Page setup: Constants and Assets sheets get an explicit portrait orientation
$constants.PageSetup.Orientation = 1
$assets.PageSetup.Orientation = 1

# Make Settings the active sheet/tab with A8 selected
$settings.Activate() | Out-Null
$settings.Range("A8").Select() | Out-Null

Write-Output "AdminInvoiceHandler config updated"
